# Updated symbol list on Fri Dec 23 22:06:52 UTC 2022 with GitHub Actions
#
# Applies the cell-value changes described by the commit diff to the
# "cryptos" worksheet: refreshed prices/ranks for most coins, the hour
# stamp (column G) moving from 21 -> 22 for every data row, a re-sort of
# three rows (KickToken / BKEXToken / CEJI) and a couple of label tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text, and whether the text looks like a
# number. Numeric-looking values are written through a "Text" number
# format first so Excel keeps the exact literal (e.g. trailing zeros like
# "0.001600") instead of silently turning them into floating point
# numbers.
$updates = @(
    @{Addr="D2"; Value="246.25"; Numeric=$true},
    @{Addr="G2"; Value="22"; Numeric=$true},
    @{Addr="D3"; Value="22.23"; Numeric=$true},
    @{Addr="G3"; Value="22"; Numeric=$true},
    @{Addr="D4"; Value="5.362"; Numeric=$true},
    @{Addr="G4"; Value="22"; Numeric=$true},
    @{Addr="D5"; Value="0.05877"; Numeric=$true},
    @{Addr="G5"; Value="22"; Numeric=$true},
    @{Addr="D6"; Value="3.392"; Numeric=$true},
    @{Addr="G6"; Value="22"; Numeric=$true},
    @{Addr="D7"; Value="6.386"; Numeric=$true},
    @{Addr="G7"; Value="22"; Numeric=$true},
    @{Addr="D8"; Value="0.8136"; Numeric=$true},
    @{Addr="G8"; Value="22"; Numeric=$true},
    @{Addr="D9"; Value="0.9540"; Numeric=$true},
    @{Addr="G9"; Value="22"; Numeric=$true},
    @{Addr="D10"; Value="0.1418"; Numeric=$true},
    @{Addr="G10"; Value="22"; Numeric=$true},
    @{Addr="D11"; Value="0.03496"; Numeric=$true},
    @{Addr="G11"; Value="22"; Numeric=$true},
    @{Addr="D12"; Value="0.07355"; Numeric=$true},
    @{Addr="G12"; Value="22"; Numeric=$true},
    @{Addr="G13"; Value="22"; Numeric=$true},
    @{Addr="G14"; Value="22"; Numeric=$true},
    @{Addr="D15"; Value="0.09401"; Numeric=$true},
    @{Addr="G15"; Value="22"; Numeric=$true},
    @{Addr="D16"; Value="0.001600"; Numeric=$true},
    @{Addr="G16"; Value="22"; Numeric=$true},
    @{Addr="D17"; Value="0.04813"; Numeric=$true},
    @{Addr="G17"; Value="22"; Numeric=$true},
    @{Addr="G18"; Value="22"; Numeric=$true},
    @{Addr="D19"; Value="0.006006"; Numeric=$true},
    @{Addr="G19"; Value="22"; Numeric=$true},
    @{Addr="D20"; Value="0.004091"; Numeric=$true},
    @{Addr="G20"; Value="22"; Numeric=$true},
    @{Addr="D21"; Value="0.0009890"; Numeric=$true},
    @{Addr="G21"; Value="22"; Numeric=$true},
    @{Addr="D22"; Value="0.00009703"; Numeric=$true},
    @{Addr="G22"; Value="22"; Numeric=$true},
    @{Addr="G23"; Value="22"; Numeric=$true},
    @{Addr="D24"; Value="2.178"; Numeric=$true},
    @{Addr="G24"; Value="22"; Numeric=$true},
    @{Addr="G25"; Value="22"; Numeric=$true},
    @{Addr="D26"; Value="0.1288"; Numeric=$true},
    @{Addr="G26"; Value="22"; Numeric=$true},
    @{Addr="D27"; Value="0.0002472"; Numeric=$true},
    @{Addr="G27"; Value="22"; Numeric=$true},
    @{Addr="G28"; Value="22"; Numeric=$true},
    @{Addr="G29"; Value="22"; Numeric=$true},
    @{Addr="G30"; Value="22"; Numeric=$true},
    @{Addr="G31"; Value="22"; Numeric=$true},
    @{Addr="G32"; Value="22"; Numeric=$true},
    @{Addr="G33"; Value="22"; Numeric=$true},
    @{Addr="G34"; Value="22"; Numeric=$true},
    @{Addr="G35"; Value="22"; Numeric=$true},
    @{Addr="G36"; Value="22"; Numeric=$true},
    @{Addr="G37"; Value="22"; Numeric=$true},
    @{Addr="G38"; Value="22"; Numeric=$true},
    @{Addr="G39"; Value="22"; Numeric=$true},
    @{Addr="D40"; Value="0.03866"; Numeric=$true},
    @{Addr="G40"; Value="22"; Numeric=$true},

    # Rows 41-43 were re-sorted: KickToken/BKEXToken/CEJI shuffled order.
    @{Addr="B41"; Value="BKEXToken"; Numeric=$false},
    @{Addr="C41"; Value="https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; Numeric=$false},
    @{Addr="D41"; Value="0.1076"; Numeric=$true},
    @{Addr="E41"; Value="40BKEXTokenBKK"; Numeric=$false},
    @{Addr="G41"; Value="22"; Numeric=$true},

    @{Addr="B42"; Value="CEJI"; Numeric=$false},
    @{Addr="C42"; Value="https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"; Numeric=$false},
    @{Addr="D42"; Value="0.002441"; Numeric=$true},
    @{Addr="E42"; Value="41CEJICEJI"; Numeric=$false},
    @{Addr="G42"; Value="22"; Numeric=$true},

    @{Addr="B43"; Value="KickToken"; Numeric=$false},
    @{Addr="C43"; Value="https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"; Numeric=$false},
    @{Addr="D43"; Value="0.003037"; Numeric=$true},
    @{Addr="E43"; Value="42KickTokenKICKWorstin24h"; Numeric=$false},
    @{Addr="G43"; Value="22"; Numeric=$true},

    @{Addr="D44"; Value="0.005786"; Numeric=$true},
    @{Addr="G44"; Value="22"; Numeric=$true},
    @{Addr="D45"; Value="0.00005663"; Numeric=$true},
    @{Addr="G45"; Value="22"; Numeric=$true},
    @{Addr="G46"; Value="22"; Numeric=$true},
    @{Addr="G47"; Value="22"; Numeric=$true},
    @{Addr="D48"; Value="0.07486"; Numeric=$true},
    @{Addr="E48"; Value="47BOLOBOLO"; Numeric=$false},
    @{Addr="G48"; Value="22"; Numeric=$true},
    @{Addr="G49"; Value="22"; Numeric=$true},
    @{Addr="G50"; Value="22"; Numeric=$true},
    @{Addr="G51"; Value="22"; Numeric=$true}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    if ($u.Numeric) {
        # Force a text number format so the literal digits (including any
        # trailing/leading zeros) are preserved verbatim instead of being
        # normalised as a floating point number.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"
